# Remove the last bullet paragraph ("Do we need to allow to delete
# company's last admin, otherwise the company cannot be deleted") from
# the "High" priority list, as part of adding the 'remember me'
# functionality item above it.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Do we need to allow to delete company") {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $d.Range($target.Range.Start, $target.Range.End)
    $r.Delete()
}
